# doc/plan.xlsx - add three more task rows to the plan table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "#100002"
$ws.Range("B4").Value = "Add mappers"

$ws.Range("A5").Value = "#100003"
$ws.Range("B5").Value = "Add simple Dao layers"

$ws.Range("A6").Value = "#100004"
$ws.Range("B6").Value = "Add Rest endpoints"

# leave the cursor where the author ended up after typing the last row
$null = $ws.Range("Q14").Select()
